# Adding the changes we made on may 9th
#
# The sensor-log sheet gains a new leading reading (the existing row 2
# values shift down into row 3, row 3's into row 4, and so on through the
# former last row, which lands in row 22), and 9 further new rows of
# accelerometer/gyroscope data are appended after that, extending the
# sheet from A1:H21 to A1:H31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: timestamp, label, ax, ay, az, gx, gy, gz -- written to rows 2..31.
$rows = @(
    @(0, "struggle", -0.01545035839080799, 0.2315296996384856, -0.5798605158925052, -0.0303905457258224, -0.02122756652534, -0.038026362657547),
    @(100, "struggle", 0.06989622116088862, -0.2621099762618542, 0.03679747134447098, 0.0122173046693205, -0.0598647929728031, -0.067195177078247),
    @(200, "struggle", 0.1602854728698729, 0.2304907962679849, -0.932947013527152, -0.1317941695451736, -0.223882108926773, -0.0047342055477201),
    @(300, "struggle", 1.401312828063962, -0.6429092288017235, -4.608263850212088, -0.5534439086914062, -2.408488988876343, -0.1953241676092147),
    @(400, "struggle", 2.035469830036161, -1.478700906038283, -5.709699153900139, -0.5908594131469727, -4.143346309661865, 0.8623890280723572),
    @(500, "struggle", 0.4589872062206247, -0.1242832839488965, -1.642685413360593, -0.3294090926647186, -3.587306261062622, 1.745547413825989),
    @(600, "struggle", -1.821219787001614, 1.541860699653622, 1.125372886657715, -1.008996725082398, -2.703842401504517, 1.318247199058533),
    @(700, "struggle", -4.537274837493896, -0.381650447845459, 0.8837652206420898, 0.1534798890352249, -0.0916297882795333, 0.2032654136419296),
    @(800, "struggle", -6.140988111495973, 2.625270426273349, 3.428681612014774, 1.313818454742432, 3.833332061767578, -1.68843150138855),
    @(900, "struggle", -6.243353843688965, 0.4006794318556786, -2.078256130218506, 0.2997821271419525, 5.671730995178223, -1.668731093406677),
    @(1000, "struggle", 1.525731027126295, -1.101333875209092, -3.408103108406068, -0.6543893814086914, 4.36127233505249, -0.7669413089752197),
    @(1100, "struggle", 2.08496618270874, -1.613636314868927, -2.943879842758179, -0.5648976564407349, 2.72186279296875, 1.022893905639648),
    @(1200, "struggle", -2.373238265514372, -0.1506620794534765, -2.15161240100861, -0.7005097270011902, 1.526399493217468, 0.2272418737411499),
    @(1300, "struggle", -3.501610040664678, 2.79044055938721, -2.906912446022055, 0.0536034256219863, -1.495856285095215, -0.2338086664676666),
    @(1400, "struggle", -4.058286607265482, 2.954759478569038, -7.582780838012708, 0.1302670091390609, -5.946773052215576, -1.086271166801453),
    @(1500, "struggle", 2.394027709960937, -1.486701488494873, 0.1956815719604492, -0.3718642294406891, -4.293160915374756, 1.783268332481384),
    @(1600, "struggle", -2.034749180078514, -1.211341693997383, -0.07141649723053026, 0.1963931769132614, -4.29804801940918, 1.745852828025818),
    @(1700, "struggle", -4.999244511127475, -0.394761115312575, -0.1581084728240968, 0.204487144947052, -2.893821477890014, 1.158658623695374),
    @(1800, "struggle", -6.522920310497289, 0.375426143407825, -0.1347707509994504, 0.1171334087848663, 0.94042706489563, -0.3061962127685547),
    @(1900, "struggle", -6.107214450836182, 0.1423146724700928, 0.9927992820739744, 0.337044894695282, 3.794236660003662, -1.184467673301697),
    @(2000, "struggle", -4.145634770393367, -0.62879066169262, 1.90212270617485, -0.3572034537792206, 5.656154155731201, -1.049466490745544),
    @(2100, "struggle", -1.618811368942268, -0.6589505374431646, 0.2403407692909387, -0.0070249503478407, 4.270253658294678, -0.0296269636601209),
    @(2200, "struggle", -0.08240008354186718, 0.4783504903316493, -3.809414207935333, 0.4137084782123565, 2.936276435852051, 0.2823724448680877),
    @(2300, "struggle", -3.95973014831543, 0.9762580394744872, -4.069071769714356, 0.0429132841527462, 1.122159481048584, 0.1867720484733581),
    @(2400, "struggle", -3.901577949523926, 1.771272063255311, -2.484678864479062, 0.06475171446800231, -1.842216849327088, -0.6108652353286743),
    @(2500, "struggle", -1.957046031951897, -0.6577051877975557, -7.9572014808655, 0.0862847194075584, -5.713422775268555, -1.346194267272949),
    @(2600, "struggle", 0.0412573218345611, -3.739429324865336, -3.584903955459609, -0.1818851232528686, -4.851491928100586, 1.392772793769836),
    @(2700, "struggle", -2.076164960861222, -2.966795355081547, -0.5832877159118695, -0.3181080818176269, -3.869678497314453, 0.9886853694915771),
    @(2800, "struggle", -5.426012933254244, -0.3007338047027568, 0.185311913490301, 0.1050688251852989, -2.216677188873291, 0.3729332387447357),
    @(2900, "struggle", -6.20224690437317, 0.8901370018720657, 1.72858691215514, 0.1996002197265625, 1.434922456741333, -0.2237294018268585)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
